# busqueda por multiples palabras
# Remove rows from the "Ejercicios" sheet that correspond to shared strings
# that should no longer exist (e.g. "acariciamiento", "acariciamientos",
# "agua", "aire", "caminar", "encuentros", "fuego", "shiva", "tierra",
# "vishnu", "Acunamiento", "Eutonia"). Rows are removed in descending
# order of row index so earlier deletions do not shift the row numbers
# of rows still pending deletion.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ejercicios")

$rowsToDelete = @(77, 73, 71, 69, 68, 28, 24, 12, 8, 7, 6, 2)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).EntireRow.Delete()
}

# Update the sheet view to match the saved selection/scroll position.
$ws.Activate()
$ws.Range("A63:XFD63").Select()
$excel.ActiveWindow.ScrollRow = 40
